$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.164.63"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.656.85"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.94"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5239"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06358"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.59"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07697"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.625"
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("D13").Value = "1.665.48"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "1.884.72"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5626"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "0.0₅8205"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.57"
$ws.Range("D18").Value = "26.151.57"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.667"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.56"
$ws.Range("E21").Value = "  +4.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "193.60"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.959"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.38"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1200"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.270"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.99"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  +2.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05493"
$ws.Range("E30").Value = "  -4.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.272"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.471"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.367"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.565"
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9536"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5692"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01591"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.873"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "1.027.99"
$ws.Range("E42").Value = "  -3.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8295"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.24"
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("D45").Value = "1.795.49"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.03"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "0.0₈104"
$ws.Range("E47").Value = "  +2.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.048"
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05198"
$ws.Range("E51").Value = "  -0.62%  "
